$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.128.38"
$ws.Range("E2").Value = "  +4.34%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.342.23"
$ws.Range("E3").Value = "  +2.81%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "415.20"
$ws.Range("E5").Value = "  +4.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "111.19"
$ws.Range("E6").Value = "  -0.14%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.587"
$ws.Range("E7").Value = "  +4.23%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.637"
$ws.Range("E9").Value = "  +2.08%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.01"
$ws.Range("E10").Value = "  +1.33%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0983"
$ws.Range("E11").Value = "  +3.62%  "

$ws.Range("E12").Value = "  +1.27%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.859.71"
$ws.Range("E13").Value = "  +2.68%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.44"
$ws.Range("E14").Value = "  +3.79%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "19.47"
$ws.Range("E15").Value = "  +1.46%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.331.34"
$ws.Range("E16").Value = "  +2.36%  "

$ws.Range("E17").Value = "  -0.71%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "58.962.72"
$ws.Range("E18").Value = "  +4.38%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.86"
$ws.Range("E19").Value = "  -1.42%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.35"
$ws.Range("E20").Value = "  +0.27%  "

$ws.Range("E21").Value = "  +4.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "13.04"
$ws.Range("E22").Value = "  -0.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "303.53"
$ws.Range("E23").Value = "  +1.19%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "75.19"
$ws.Range("E24").Value = "  -0.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.27"
$ws.Range("E25").Value = "  +0.79%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "28.73"
$ws.Range("E26").Value = "  +1.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.93"
$ws.Range("E27").Value = "  -3.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.35"
$ws.Range("E28").Value = "  -0.23%  "

$ws.Range("E29").Value = "  -0.02%  "

$ws.Range("E30").Value = "  +0.07%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.113"
$ws.Range("E31").Value = "  +1.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.46"
$ws.Range("E32").Value = "  +2.74%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "40.38"
$ws.Range("E33").Value = "  +9.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0540"
$ws.Range("E34").Value = "  +10.28%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.13"
$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "52.14"
$ws.Range("E36").Value = "  +1.17%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.22"
$ws.Range("E37").Value = "  +4.05%  "

$ws.Range("E38").Value = "  +0.00%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.50"
$ws.Range("E39").Value = "  -1.25%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "137.76"
$ws.Range("E40").Value = "  +2.07%  "

$ws.Range("E41").Value = "  +1.59%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.92"
$ws.Range("E42").Value = "  -0.89%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.97"
$ws.Range("E43").Value = "  -1.15%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "16.82"
$ws.Range("E44").Value = "  -4.95%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.279"
$ws.Range("E45").Value = "  -1.63%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.26"
$ws.Range("E46").Value = "  +8.07%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "22.47"
$ws.Range("E47").Value = "  +0.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.197.81"
$ws.Range("E48").Value = "  +2.49%  "

$ws.Range("E49").Value = "  -0.55%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.95"
$ws.Range("E50").Value = "  -10.41%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.39"
$ws.Range("E51").Value = "  +5.69%  "
